# Auto-generated edit script: updates crypto price/volume figures
# (and two name/link swaps) to match the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.355.17"
$ws.Range("E2").Value = "  +2.84%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.96"
$ws.Range("E3").Value = "  +2.81%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.13"
$ws.Range("E5").Value = "  +1.35%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.02"
$ws.Range("E6").Value = "  +5.47%  "
# Row 7
$ws.Range("E7").Value = "  -0.11%  "
# Row 8
$ws.Range("E8").Value = "  +0.76%  "
# Row 9
$ws.Range("E9").Value = "  +0.57%  "
# Row 10
$ws.Range("E10").Value = "  +3.19%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.398"
$ws.Range("E11").Value = "  +2.75%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.085.11"
$ws.Range("E12").Value = "  +2.85%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.86"
$ws.Range("E13").Value = "  +4.24%  "
# Row 14
$ws.Range("E14").Value = "  -0.44%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.498.60"
$ws.Range("E15").Value = "  +3.23%  "
# Row 16
$ws.Range("E16").Value = "  +2.46%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.371.65"
$ws.Range("E17").Value = "  +2.83%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.31"
$ws.Range("E18").Value = "  +2.61%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.40"
$ws.Range("E19").Value = "  +5.60%  "
# Row 20
$ws.Range("E20").Value = "  +4.67%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.35"
$ws.Range("E21").Value = "  +0.50%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.566"
$ws.Range("E22").Value = "  +1.72%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.24"
$ws.Range("E23").Value = "  -0.16%  "
# Row 24
$ws.Range("E24").Value = "  -0.09%  "
# Row 25
$ws.Range("E25").Value = "  +5.80%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.629.17"
$ws.Range("E26").Value = "  +3.16%  "
# Row 27
$ws.Range("E27").Value = "  -4.23%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.87"
$ws.Range("E28").Value = "  +8.65%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.29%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.31"
$ws.Range("E30").Value = "  +3.34%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.48"
$ws.Range("E31").Value = "  +8.14%  "
# Row 32
$ws.Range("E32").Value = "  +0.41%  "
# Row 33
$ws.Range("E33").Value = "  -0.05%  "
# Row 34
$ws.Range("E34").Value = "  +2.00%  "
# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  +7.09%  "
# Row 36
$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.65"
$ws.Range("E36").Value = "  +26.97%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.15"
$ws.Range("E37").Value = "  +3.05%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "171.39"
$ws.Range("E38").Value = "  +2.21%  "
# Row 39
$ws.Range("E39").Value = "  +7.46%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.524.89"
$ws.Range("E40").Value = "  +2.79%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0772"
$ws.Range("E41").Value = "  +0.66%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.809"
$ws.Range("E42").Value = "  +3.77%  "
# Row 43
$ws.Range("E43").Value = "  +2.06%  "
# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("E44").Value = "  +4.50%  "
# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.46"
$ws.Range("E45").Value = "  +0.10%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.22"
$ws.Range("E46").Value = "  +7.90%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.619.44"
$ws.Range("E47").Value = "  +6.48%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.75"
$ws.Range("E48").Value = "  +5.40%  "
# Row 49
$ws.Range("E49").Value = "  +13.68%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("E50").Value = "  +1.26%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0272"
$ws.Range("E51").Value = "  +3.65%  "
